# Auto-generated edit script based on target diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("C2").Value = 46063

# --- Row 3 ---
$ws.Range("C3").Value = 46063

# --- Row 4 ---
$ws.Range("A4").Value = 'A 1782-2024'
$ws.Range("B4").Value = 45307
$ws.Range("C4").Value = 46063
$ws.Range("G4").Value = 2.7
$ws.Range("S4").Formula = 'HYPERLINK("https://klasma.github.io/Logging_1278/artfynd/A 1782-2024 artfynd.xlsx", "A 1782-2024")'
$ws.Range("T4").Formula = 'HYPERLINK("https://klasma.github.io/Logging_1278/kartor/A 1782-2024 karta.png", "A 1782-2024")'
$ws.Range("V4").Formula = 'HYPERLINK("https://klasma.github.io/Logging_1278/klagomål/A 1782-2024 FSC-klagomål.docx", "A 1782-2024")'
$ws.Range("W4").Formula = 'HYPERLINK("https://klasma.github.io/Logging_1278/klagomålsmail/A 1782-2024 FSC-klagomål mail.docx", "A 1782-2024")'
$ws.Range("X4").Formula = 'HYPERLINK("https://klasma.github.io/Logging_1278/tillsyn/A 1782-2024 tillsynsbegäran.docx", "A 1782-2024")'
$ws.Range("Y4").Formula = 'HYPERLINK("https://klasma.github.io/Logging_1278/tillsynsmail/A 1782-2024 tillsynsbegäran mail.docx", "A 1782-2024")'

# --- Row 5 ---
$ws.Range("C5").Value = 46063

# --- Row 6 ---
$ws.Range("A6").Value = 'A 389-2023'
$ws.Range("B6").Value = 44929
$ws.Range("C6").Value = 46063
$ws.Range("G6").Value = 2.5
$ws.Range("S6").Formula = 'HYPERLINK("https://klasma.github.io/Logging_1278/artfynd/A 389-2023 artfynd.xlsx", "A 389-2023")'
$ws.Range("T6").Formula = 'HYPERLINK("https://klasma.github.io/Logging_1278/kartor/A 389-2023 karta.png", "A 389-2023")'
$ws.Range("V6").Formula = 'HYPERLINK("https://klasma.github.io/Logging_1278/klagomål/A 389-2023 FSC-klagomål.docx", "A 389-2023")'
$ws.Range("W6").Formula = 'HYPERLINK("https://klasma.github.io/Logging_1278/klagomålsmail/A 389-2023 FSC-klagomål mail.docx", "A 389-2023")'
$ws.Range("X6").Formula = 'HYPERLINK("https://klasma.github.io/Logging_1278/tillsyn/A 389-2023 tillsynsbegäran.docx", "A 389-2023")'
$ws.Range("Y6").Formula = 'HYPERLINK("https://klasma.github.io/Logging_1278/tillsynsmail/A 389-2023 tillsynsbegäran mail.docx", "A 389-2023")'

# --- Row 7 ---
$ws.Range("C7").Value = 46063

# --- Row 8 ---
$ws.Range("A8").Value = 'A 48974-2023'
$ws.Range("B8").Value = 45209
$ws.Range("C8").Value = 46063
$ws.Range("G8").Value = 4.5

# --- Row 9 ---
$ws.Range("A9").Value = 'A 4822-2023'
$ws.Range("B9").Value = 44957
$ws.Range("C9").Value = 46063
$ws.Range("G9").Value = 2.2

# --- Row 10 ---
$ws.Range("A10").Value = 'A 1531-2022'
$ws.Range("B10").Value = 44573
$ws.Range("C10").Value = 46063
$ws.Range("G10").Value = 1.6

# --- Row 11 ---
$ws.Range("A11").Value = 'A 4486-2024'
$ws.Range("B11").Value = 45327
$ws.Range("C11").Value = 46063
$ws.Range("G11").Value = 0.6

# --- Row 12 ---
$ws.Range("A12").Value = 'A 32610-2024'
$ws.Range("B12").Value = 45513
$ws.Range("C12").Value = 46063
$ws.Range("F12").ClearContents()
$ws.Range("G12").Value = 0.5

# --- Row 13 ---
$ws.Range("A13").Value = 'A 10710-2025'
$ws.Range("B13").Value = 45722
$ws.Range("C13").Value = 46063
$ws.Range("F13").Value = 'Kommuner'
$ws.Range("G13").Value = 1.8

# --- Row 14 ---
$ws.Range("A14").Value = 'A 635-2023'
$ws.Range("B14").Value = 44930
$ws.Range("C14").Value = 46063
$ws.Range("G14").Value = 0.5

# --- Row 15 ---
$ws.Range("A15").Value = 'A 48181-2024'
$ws.Range("B15").Value = 45589
$ws.Range("C15").Value = 46063
$ws.Range("G15").Value = 0.7

# --- Row 16 ---
$ws.Range("A16").Value = 'A 34400-2025'
$ws.Range("B16").Value = 45846.61351851852
$ws.Range("C16").Value = 46063
$ws.Range("G16").Value = 1.3

# --- Row 17 ---
$ws.Range("A17").Value = 'A 34401-2025'
$ws.Range("B17").Value = 45846.6140162037
$ws.Range("C17").Value = 46063
$ws.Range("G17").Value = 2.8

# --- Row 18 ---
$ws.Range("A18").Value = 'A 18328-2025'
$ws.Range("B18").Value = 45762
$ws.Range("C18").Value = 46063
$ws.Range("G18").Value = 1.8

# --- Row 19 ---
$ws.Range("A19").Value = 'A 5817-2025'
$ws.Range("B19").Value = 45694.74113425926
$ws.Range("C19").Value = 46063
$ws.Range("G19").Value = 1.2

# --- Row 20 ---
$ws.Range("A20").Value = 'A 24-2023'
$ws.Range("B20").Value = 44928
$ws.Range("C20").Value = 46063
$ws.Range("G20").Value = 0.5

# --- Row 21 ---
$ws.Range("C21").Value = 46063

# --- Row 22 ---
$ws.Range("A22").Value = 'A 28260-2023'
$ws.Range("B22").Value = 45099
$ws.Range("C22").Value = 46063
$ws.Range("G22").Value = 5

# --- Row 23 ---
$ws.Range("A23").Value = 'A 4481-2024'
$ws.Range("B23").Value = 45327
$ws.Range("C23").Value = 46063
$ws.Range("G23").Value = 1

# --- Row 24 ---
$ws.Range("A24").Value = 'A 18332-2025'
$ws.Range("B24").Value = 45762
$ws.Range("C24").Value = 46063
$ws.Range("G24").Value = 2.5

# --- Row 25 ---
$ws.Range("A25").Value = 'A 4256-2025'
$ws.Range("B25").Value = 45685
$ws.Range("C25").Value = 46063
$ws.Range("G25").Value = 2

# --- Row 26 ---
$ws.Range("A26").Value = 'A 11517-2024'
$ws.Range("B26").Value = 45372
$ws.Range("C26").Value = 46063

# --- Row 27 ---
$ws.Range("A27").Value = 'A 18434-2023'
$ws.Range("B27").Value = 45042
$ws.Range("C27").Value = 46063
$ws.Range("G27").Value = 0.7

# --- Row 28 ---
$ws.Range("A28").Value = 'A 7731-2026'
$ws.Range("B28").Value = 46062.52008101852
$ws.Range("C28").Value = 46063
$ws.Range("G28").Value = 5.9
$ws.Range("R28").ClearContents()
$ws.Range("R28").WrapText = $true

# --- Row 29 ---
$ws.Range("A29").Value = 'A 42412-2025'
$ws.Range("B29").Value = 45905
$ws.Range("C29").Value = 46063
$ws.Range("G29").Value = 3.1

# --- Row 30 ---
$ws.Range("A30").Value = 'A 42417-2025'
$ws.Range("B30").Value = 45905
$ws.Range("C30").Value = 46063
$ws.Range("G30").Value = 1.1

# --- Row 31 ---
$ws.Range("A31").Value = 'A 7727-2026'
$ws.Range("B31").Value = 46062.50420138889
$ws.Range("C31").Value = 46063
$ws.Range("G31").Value = 1.9
$ws.Range("R31").ClearContents()
$ws.Range("R31").WrapText = $true

# --- Row 32 ---
$ws.Range("A32").Value = 'A 53131-2021'
$ws.Range("B32").Value = 44468
$ws.Range("C32").Value = 46063
$ws.Range("G32").Value = 1.3

# --- Row 33 ---
$ws.Range("A33").Value = 'A 4487-2024'
$ws.Range("C33").Value = 46063
$ws.Range("G33").Value = 1.9

# --- Row 34 ---
$ws.Range("A34").Value = 'A 18327-2025'
$ws.Range("B34").Value = 45762
$ws.Range("B34").NumberFormat = "YYYY-MM-DD"
$ws.Range("C34").Value = 46063
$ws.Range("C34").NumberFormat = "YYYY-MM-DD"
$ws.Range("D34").Value = 'SKÅNE LÄN'
$ws.Range("E34").Value = 'BÅSTAD'
$ws.Range("G34").Value = 0.6
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = 0
$ws.Range("N34").Value = 0
$ws.Range("O34").Value = 0
$ws.Range("P34").Value = 0
$ws.Range("Q34").Value = 0
$ws.Range("R34").ClearContents()
$ws.Range("R34").WrapText = $true

# --- Row 35 ---
$ws.Range("A35").Value = 'A 4493-2024'
$ws.Range("B35").Value = 45327
$ws.Range("B35").NumberFormat = "YYYY-MM-DD"
$ws.Range("C35").Value = 46063
$ws.Range("C35").NumberFormat = "YYYY-MM-DD"
$ws.Range("D35").Value = 'SKÅNE LÄN'
$ws.Range("E35").Value = 'BÅSTAD'
$ws.Range("G35").Value = 1.8
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = 0
$ws.Range("N35").Value = 0
$ws.Range("O35").Value = 0
$ws.Range("P35").Value = 0
$ws.Range("Q35").Value = 0
$ws.Range("R35").ClearContents()
$ws.Range("R35").WrapText = $true

# --- Row height adjustments ---
$ws.Rows(33).RowHeight = 15
$ws.Rows(34).RowHeight = 15
